$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.735.05'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '2.369.85'
$ws.Range("E3").Value = '  +1.58%  '
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -9.12%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.619'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0921'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("E12").Value = '  -5.90%  '
$ws.Range("E13").Value = '  -5.87%  '
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.27'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.50%  '
$ws.Range("D16").Value = '2.733.50'
$ws.Range("E16").Value = '  +2.29%  '
$ws.Range("D17").Value = '2.370.18'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = '42.725.90'
$ws.Range("E18").Value = '  -0.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.08%  '
$ws.Range("E20").Value = '  -2.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.97%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '271.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +6.73%  '
$ws.Range("E24").Value = '  -7.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.76'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.92%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.20'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '172.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.80%  '
$ws.Range("E31").Value = '  -1.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0898'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -10.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.58'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0357'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.87'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.94%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.85'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.47%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.104'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.37%  '
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.14'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +50.75%  '
$ws.Range("E43").Value = '  -3.59%  '
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '68.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.12%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '115.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.81'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.86%  '
$ws.Range("E49").Value = '  -3.14%  '
$ws.Range("D50").Value = '1.602.07'
$ws.Range("E50").Value = '  +7.90%  '
$ws.Range("E51").Value = '  -3.13%  '
